$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the data for rows 6-9 (analysis #05 .. #08 / cpx #5 .. #8 block),
# leaving the empty, still-styled cells behind (matches the diff: cells keep
# their "s" style attribute but lose their t="s"/<v> content).
$ws.Range("A6:P9").Value = ""

# Update the selection to match the edited sheet view.
$ws.Range("A6:Q11").Select()
